# Applies the "Added a graph. Need to fix x axis." edit to Sheet1.
# (Per the captured OOXML diff, the concrete cell-level changes are a
# relabeling/restructuring of the input section of the sheet: the model
# gets English labels, a "Loan & Capital:" sub-header, a "Down to:" /
# "Amortization amount per year" / "Add level" header row, and an
# "Interest deduction" percentage input.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: section header -------------------------------------------------
$ws.Range("A6").Value = "Amortizing model"

# --- Row 7: old "Amorterings modell" label is gone; row becomes the new
#            column header row for the level table -------------------------
$ws.Range("A7").ClearContents()
$ws.Range("B7").Value = "Down to:"
$ws.Range("C7").Value = "Amortization amount per year"
$ws.Range("E7").Value = "Add level"

# --- Rows 8-9: level rows, labels translated, numeric inputs unchanged -----
$ws.Range("A8").Value = "Level 1"
$ws.Range("A9").Value = "Level 2"

# --- Row 10: "Möjlighet att lägga till fler" removed ------------------------
$ws.Range("A10").ClearContents()

# --- Row 12: new "Loan & Capital:" sub-header -------------------------------
$ws.Range("A12").Value = "Loan & Capital:"

# --- Row 17: interest-deduction label + new percentage input ---------------
$ws.Range("A17").Value = "Interest deduction"
$ws.Range("B17").Value = 0.3
$ws.Range("B17").NumberFormat = "0%"

# --- Row 34: swap which label sits in Output vs "Eventuell output" section -
$ws.Range("A34").Value = "Eventuell output:"
$ws.Range("B34").Value = "Datagrid med måndskostnad för konfigurerbart intervall"

# --- Column widths / layout --------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 23

# --- View state: scrolled down a bit, selection parked on A18 --------------
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("A18").Select()
